$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($i = 2; $i -le $lastRow; $i++) {
    $typeCell = $ws.Cells.Item($i, 2)
    $valueCell = $ws.Cells.Item($i, 3)

    $typeVal = $typeCell.Value2
    $sportVal = $valueCell.Value2

    if ($typeVal -eq $null -or $sportVal -eq $null) {
        continue
    }

    # e.g. "club-sports" -> "club", "uil-sports" -> "uil"
    $category = $typeVal.Split("-")[0]

    # e.g. "Football-Boys" -> "Boys" -> "boys"
    $sportParts = $sportVal.Split("-")
    $sportName = $sportParts[0]
    $gender = $sportParts[1].ToLower()

    $newType = "sports_" + $category + "_" + $gender

    $typeCell.Value = $newType
    $valueCell.Value = $sportName
}
